# Add two new data-year columns (R = 2021, S = 2022) to the "3.9.1.1" table
# on sheet1, mirroring the layout/formatting already used by column Q (2020).
# For every data row we copy column Q's cell (to inherit its number format /
# font / alignment) into the new R and S cells of that row, then overwrite
# the copied value with the correct figure from the source dataset. Cells
# that hold the "no data" placeholder get the literal text "-" (which Excel
# maps back onto the existing shared string used for every other "-" cell
# in the sheet) instead of a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - year headers
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021
$ws.Range("Q4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 2022

# Row 5
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 0.8
$ws.Range("Q5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 0.5

# Row 6
$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("R6").Value = 0.4
$ws.Range("Q6").Copy($ws.Range("S6"))
$ws.Range("S6").Value = 0.2

# Row 7
$ws.Range("Q7").Copy($ws.Range("R7"))
$ws.Range("R7").Value = 1.2
$ws.Range("Q7").Copy($ws.Range("S7"))
$ws.Range("S7").Value = 0.7

# Row 8
$ws.Range("Q8").Copy($ws.Range("R8"))
$ws.Range("R8").Value = 0.2
$ws.Range("Q8").Copy($ws.Range("S8"))
$ws.Range("S8").Value = 0.2

# Row 9
$ws.Range("Q9").Copy($ws.Range("R9"))
$ws.Range("R9").Value = "-"
$ws.Range("Q9").Copy($ws.Range("S9"))
$ws.Range("S9").Value = "-"

# Row 10
$ws.Range("Q10").Copy($ws.Range("R10"))
$ws.Range("R10").Value = 0.4
$ws.Range("Q10").Copy($ws.Range("S10"))
$ws.Range("S10").Value = 0.4

# Row 11
$ws.Range("Q11").Copy($ws.Range("R11"))
$ws.Range("R11").Value = 0.6
$ws.Range("Q11").Copy($ws.Range("S11"))
$ws.Range("S11").Value = 0.5

# Row 12
$ws.Range("Q12").Copy($ws.Range("R12"))
$ws.Range("R12").Value = 0.8
$ws.Range("Q12").Copy($ws.Range("S12"))
$ws.Range("S12").Value = 0.3

# Row 13
$ws.Range("Q13").Copy($ws.Range("R13"))
$ws.Range("R13").Value = 0.5
$ws.Range("Q13").Copy($ws.Range("S13"))
$ws.Range("S13").Value = 0.6

# Row 14
$ws.Range("Q14").Copy($ws.Range("R14"))
$ws.Range("R14").Value = 0.4
$ws.Range("Q14").Copy($ws.Range("S14"))
$ws.Range("S14").Value = 0.7

# Row 15
$ws.Range("Q15").Copy($ws.Range("R15"))
$ws.Range("R15").Value = "-"
$ws.Range("Q15").Copy($ws.Range("S15"))
$ws.Range("S15").Value = 0.4

# Row 16
$ws.Range("Q16").Copy($ws.Range("R16"))
$ws.Range("R16").Value = 0.8
$ws.Range("Q16").Copy($ws.Range("S16"))
$ws.Range("S16").Value = 1.1000000000000001

# Row 17
$ws.Range("Q17").Copy($ws.Range("R17"))
$ws.Range("R17").Value = 0.3
$ws.Range("Q17").Copy($ws.Range("S17"))
$ws.Range("S17").Value = "-"

# Row 18
$ws.Range("Q18").Copy($ws.Range("R18"))
$ws.Range("R18").Value = 0.7
$ws.Range("Q18").Copy($ws.Range("S18"))
$ws.Range("S18").Value = "-"

# Row 19
$ws.Range("Q19").Copy($ws.Range("R19"))
$ws.Range("R19").Value = "-"
$ws.Range("Q19").Copy($ws.Range("S19"))
$ws.Range("S19").Value = "-"

# Row 20
$ws.Range("Q20").Copy($ws.Range("R20"))
$ws.Range("R20").Value = 0.5
$ws.Range("Q20").Copy($ws.Range("S20"))
$ws.Range("S20").Value = 0.4

# Row 21
$ws.Range("Q21").Copy($ws.Range("R21"))
$ws.Range("R21").Value = 0.1
$ws.Range("Q21").Copy($ws.Range("S21"))
$ws.Range("S21").Value = 0.4

# Row 22
$ws.Range("Q22").Copy($ws.Range("R22"))
$ws.Range("R22").Value = 0.8
$ws.Range("Q22").Copy($ws.Range("S22"))
$ws.Range("S22").Value = 0.4

# Row 23
$ws.Range("Q23").Copy($ws.Range("R23"))
$ws.Range("R23").Value = 1.1000000000000001
$ws.Range("Q23").Copy($ws.Range("S23"))
$ws.Range("S23").Value = 0.4

# Row 24 (S24 uses the formatting of the "-" cells in this block, same as
# the original workbook, rather than Q24's own style)
$ws.Range("Q24").Copy($ws.Range("R24"))
$ws.Range("R24").Value = 1.5
$ws.Range("Q25").Copy($ws.Range("S24"))
$ws.Range("S24").Value = "-"

# Row 25
$ws.Range("Q25").Copy($ws.Range("R25"))
$ws.Range("R25").Value = 0.7
$ws.Range("Q25").Copy($ws.Range("S25"))
$ws.Range("S25").Value = 0.7

# Row 26
$ws.Range("Q26").Copy($ws.Range("R26"))
$ws.Range("R26").Value = 2.2000000000000002
$ws.Range("Q26").Copy($ws.Range("S26"))
$ws.Range("S26").Value = 1

# Row 27
$ws.Range("Q27").Copy($ws.Range("R27"))
$ws.Range("R27").Value = 1
$ws.Range("Q27").Copy($ws.Range("S27"))
$ws.Range("S27").Value = 0.4

# Row 28
$ws.Range("Q28").Copy($ws.Range("R28"))
$ws.Range("R28").Value = 3.5
$ws.Range("Q28").Copy($ws.Range("S28"))
$ws.Range("S28").Value = 1.7

# Row 29
$ws.Range("Q29").Copy($ws.Range("R29"))
$ws.Range("R29").Value = 0.8
$ws.Range("Q29").Copy($ws.Range("S29"))
$ws.Range("S29").Value = 0.3

# Row 30
$ws.Range("Q30").Copy($ws.Range("R30"))
$ws.Range("R30").Value = 0.2
$ws.Range("Q30").Copy($ws.Range("S30"))
$ws.Range("S30").Value = 0

# Row 31
$ws.Range("Q31").Copy($ws.Range("R31"))
$ws.Range("R31").Value = 1.6
$ws.Range("Q31").Copy($ws.Range("S31"))
$ws.Range("S31").Value = 0.6

# Row 32
$ws.Range("Q32").Copy($ws.Range("R32"))
$ws.Range("R32").Value = 0.3
$ws.Range("Q32").Copy($ws.Range("S32"))
$ws.Range("S32").Value = "-"

# Row 33
$ws.Range("Q33").Copy($ws.Range("R33"))
$ws.Range("R33").Value = "-"
$ws.Range("Q33").Copy($ws.Range("S33"))
$ws.Range("S33").Value = "-"

# Row 34
$ws.Range("Q34").Copy($ws.Range("R34"))
$ws.Range("R34").Value = 0.6
$ws.Range("Q34").Copy($ws.Range("S34"))
$ws.Range("S34").Value = "-"

# Match the saved selection recorded in the workbook (cell T6 active).
$ws.Activate()
$ws.Range("T6").Select()
